# Applies the BGR model update: re-ranking of cost classes for a handful of
# solar/wind resource rows (names, descriptions and the lcoe_class / P column)
# on the "solar" and "wind" worksheets.

function Set-Row {
    param($ws, $row, $resource, $kind, $costClass)

    $name = "e_" + $resource + "_c" + $costClass
    $desc = $kind + " resource -- CF class " + $resource + " -- cost class " + $costClass

    $ws.Range("C$row").Value = $name
    $ws.Range("K$row").Value = $name
    $ws.Range("D$row").Value = $desc
    $ws.Range("P$row").Value = $costClass
}

$wb = $excel.ActiveWorkbook

# --- solar sheet ---
$solar = $wb.Worksheets.Item("solar")
Set-Row $solar 4 "spv-BGR_16" "solar" 4
Set-Row $solar 5 "spv-BGR_16" "solar" 2
Set-Row $solar 6 "spv-BGR_16" "solar" 3

# --- wind sheet ---
$wind = $wb.Worksheets.Item("wind")
Set-Row $wind 4  "won-BGR_29" "wind" 2
Set-Row $wind 5  "won-BGR_29" "wind" 3

Set-Row $wind 16 "won-BGR_25" "wind" 1
Set-Row $wind 17 "won-BGR_25" "wind" 2

Set-Row $wind 18 "won-BGR_24" "wind" 2
Set-Row $wind 19 "won-BGR_24" "wind" 3
Set-Row $wind 20 "won-BGR_24" "wind" 1

Set-Row $wind 47 "won-BGR_17" "wind" 1
Set-Row $wind 48 "won-BGR_17" "wind" 2

$wb.Save()
